$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "58.097.25"
$ws.Range("E2").Value2 = "  -3.75%  "

# Row 3
$ws.Range("D3").Value2 = "3.130.35"
$ws.Range("E3").Value2 = "  -5.20%  "

# Row 4
$ws.Range("E4").Value2 = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "523.35"
$ws.Range("E5").Value2 = "  -6.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "134.16"
$ws.Range("E6").Value2 = "  -5.57%  "

# Row 7
$ws.Range("E7").Value2 = "  -0.08%  "

# Row 8
$ws.Range("D8").Value2 = "3.130.07"
$ws.Range("E8").Value2 = "  -5.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.438"
$ws.Range("E9").Value2 = "  -6.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.18"
$ws.Range("E10").Value2 = "  -8.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.108"
$ws.Range("E11").Value2 = "  -9.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.376"
$ws.Range("E12").Value2 = "  -7.69%  "

# Row 13
$ws.Range("D13").Value2 = "3.666.97"
$ws.Range("E13").Value2 = "  -5.17%  "

# Row 14
$ws.Range("E14").Value2 = "  -1.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "25.43"
$ws.Range("E15").Value2 = "  -5.41%  "

# Row 16
$ws.Range("D16").Value2 = "3.130.33"
$ws.Range("E16").Value2 = "  -5.20%  "

# Row 17
$ws.Range("D17").Value2 = "58.090.89"
$ws.Range("E17").Value2 = "  -3.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.0000151"
$ws.Range("E18").Value2 = "  -8.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "5.77"
$ws.Range("E19").Value2 = "  -5.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "13.02"
$ws.Range("E20").Value2 = "  -7.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.89"
$ws.Range("E21").Value2 = "  -8.72%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "342.00"
$ws.Range("E22").Value2 = "  -8.79%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.999"
$ws.Range("E23").Value2 = "  -0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "67.73"
$ws.Range("E24").Value2 = "  -8.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.504"
$ws.Range("E25").Value2 = "  -5.86%  "

# Row 26
$ws.Range("D26").Value2 = "3.260.82"
$ws.Range("E26").Value2 = "  -5.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.166"
$ws.Range("E27").Value2 = "  -3.40%  "

# Row 28
$ws.Range("B28").Value2 = "Binance-PegBSC-USD"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.00"
$ws.Range("E28").Value2 = "  +0.76%  "

# Row 29
$ws.Range("B29").Value2 = "PEPE"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value2 = "0.0₃0946"
$ws.Range("E29").Value2 = "  -7.82%  "

# Row 30
$ws.Range("B30").Value2 = "USDe"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.999"
$ws.Range("E30").Value2 = "  -0.01%  "

# Row 31
$ws.Range("B31").Value2 = "RenderToken"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.74"
$ws.Range("E31").Value2 = "  -6.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.86"
$ws.Range("E32").Value2 = "  -8.89%  "

# Row 33
$ws.Range("B33").Value2 = "EthereumClassic"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "21.38"
$ws.Range("E33").Value2 = "  -5.44%  "

# Row 34
$ws.Range("B34").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "6.77"
$ws.Range("E34").Value2 = "  -10.43%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.20"
$ws.Range("E35").Value2 = "  -3.57%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "4.77"
$ws.Range("E36").Value2 = "  -6.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "156.77"
$ws.Range("E37").Value2 = "  -5.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "6.18"
$ws.Range("E38").Value2 = "  -7.55%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.36"
$ws.Range("E39").Value2 = "  -11.32%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.0686"
$ws.Range("E40").Value2 = "  -5.98%  "

# Row 41
$ws.Range("B41").Value2 = "RenzoRestakedETH"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value2 = "3.161.28"
$ws.Range("E41").Value2 = "  -5.15%  "

# Row 42
$ws.Range("B42").Value2 = "OKB"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "40.24"
$ws.Range("E42").Value2 = "  -4.06%  "

# Row 43
$ws.Range("B43").Value2 = "EnergySwap"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "23.91"
$ws.Range("E43").Value2 = "  -10.16%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.688"
$ws.Range("E44").Value2 = "  -8.31%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.07"
$ws.Range("E45").Value2 = "  -3.50%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.88"
$ws.Range("E46").Value2 = "  -6.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.999"
$ws.Range("E47").Value2 = "  -0.03%  "

# Row 48
$ws.Range("B48").Value2 = "Maker"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value2 = "2.257.06"
$ws.Range("E48").Value2 = "  -3.92%  "

# Row 49
$ws.Range("B49").Value2 = "Stacks"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.42"
$ws.Range("E49").Value2 = "  -9.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "6.17"
$ws.Range("E50").Value2 = "  -3.50%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "20.62"
$ws.Range("E51").Value2 = "  -3.24%  "
